# Rows 3-6 on the "Artfynd" sheet got re-sorted/re-assigned: the data that
# used to live in row 5 now lives in row 3, row 3's old data moved to row 4,
# row 6's old data moved to row 5, and row 4's old data moved to row 6.
# Only columns A, B, D, E, F, G, H, Q, R change; everything else in those
# rows (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AS, AT, AW, AX, AY)
# is identical across the four rows already, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that actually change, keyed
# by their current row, so the writes below don't clobber data we still
# need to read later.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$before = @{}
foreach ($r in 3..6) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowVals
}

# new row -> old row it takes its data from
$mapping = @{ 3 = 5; 4 = 3; 5 = 6; 6 = 4 }

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $before[$oldRow][$col]
    }
}
